# Apply "test P7 with -10 percent" edits across the result sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "general": summary KPIs -----------------------------------
$ws = $wb.Worksheets.Item("general")
$ws.Range("B3").Value  = 441.4078650152845      # objValue
$ws.Range("B4").Value  = 0.0130000114440918     # runtime
$ws.Range("B6").Value  = 45.77786501528455      # Z1
$ws.Range("B10").Value = 395.63                 # Z5

# --- Sheet "x": assignment indices ------------------------------------
$ws = $wb.Worksheets.Item("x")
$ws.Range("B5").Value  = 11
$ws.Range("B8").Value  = 10
$ws.Range("B11").Value = 13
$ws.Range("B12").Value = 12
$ws.Range("B13").Value = 7
$ws.Range("B14").Value = 1

# --- Sheet "U" ---------------------------------------------------------
$ws = $wb.Worksheets.Item("U")
$ws.Range("B3").Value = 2

# --- Sheet "TBar" --------------------------------------------------------
$ws = $wb.Worksheets.Item("TBar")
$ws.Range("B3").Value  = 22.61068956408072
$ws.Range("B4").Value  = 10
$ws.Range("B9").Value  = 23.22876137241512
$ws.Range("B11").Value = 20
$ws.Range("B12").Value = 20
$ws.Range("B13").Value = 30
$ws.Range("B14").Value = 25.21630137166173
$ws.Range("B15").Value = 28.21630585843227

# --- Sheet "Q" -----------------------------------------------------------
$ws = $wb.Worksheets.Item("Q")
$ws.Range("C7").Value  = 297.2549999999997
$ws.Range("C8").Value  = 314.6649999999998
$ws.Range("C9").Value  = 311.2049999999998
$ws.Range("C10").Value = 316.3799999999999
$ws.Range("C11").Value = 301.5249999999998
$ws.Range("C12").Value = 67.77500000000072
$ws.Range("C13").Value = 73.77000000000072
$ws.Range("C14").Value = 74.03500000000074
$ws.Range("C15").Value = 73.04500000000073
$ws.Range("C16").Value = 72.66500000000073
$ws.Range("C22").Value = 129.3599999999996
$ws.Range("C25").Value = 136.5699999999998
$ws.Range("C37").Value = 236.3350000000021
$ws.Range("C38").Value = 246.4550000000021
$ws.Range("C39").Value = 231.7250000000021
$ws.Range("C40").Value = 253.5450000000021
$ws.Range("C41").Value = 239.25
$ws.Range("C47").Value = 153.2600000000012
$ws.Range("C48").Value = 161.7350000000012
$ws.Range("C49").Value = 153.75
$ws.Range("C50").Value = 163.7750000000012
$ws.Range("C51").Value = 157.3950000000012
$ws.Range("C52").Value = 187.3099999999989
$ws.Range("C53").Value = 197.074999999999
$ws.Range("C54").Value = 197.934999999999
$ws.Range("C55").Value = 197.2249999999991
$ws.Range("C56").Value = 185.719999999999
$ws.Range("C57").Value = 297.2549999999997
$ws.Range("C58").Value = 314.6649999999998
$ws.Range("C59").Value = 311.2049999999998
$ws.Range("C60").Value = 316.3799999999999
$ws.Range("C61").Value = 301.5249999999998
$ws.Range("C62").Value = 212.0549999999987
$ws.Range("C63").Value = 215.8299999999987
$ws.Range("C64").Value = 177.0399999999987
$ws.Range("C65").Value = 198
$ws.Range("C66").Value = 184.7
$ws.Range("C67").Value = 236.3350000000021
$ws.Range("C68").Value = 246.4550000000021
$ws.Range("C69").Value = 231.7250000000021
$ws.Range("C70").Value = 253.5450000000021
$ws.Range("C71").Value = 239.25

# --- Sheet "L" -----------------------------------------------------------
$ws = $wb.Worksheets.Item("L")
$ws.Range("C7").Value  = 13.7
$ws.Range("C8").Value  = 6.91
$ws.Range("C9").Value  = 10.68
$ws.Range("C10").Value = 7.39
$ws.Range("C11").Value = 14.68
